$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute("Gwiazdozbiór Wolarza", $true, $false, $false, $false, $false,
              $true, 1, $false, "Gwiazdozbiór Butów", 2)
